$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = "Datos actualizados a 21 de Mayo de 2020 a las 20:05"

# Row 4: Estados Unidos
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 1603926
$ws.Cells.Item(4, 3).Value = 11203
$ws.Cells.Item(4, 4).Value = 371795
$ws.Cells.Item(4, 5).Value = 1136630
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 565
$ws.Cells.Item(4, 8).Value = 95501

# Row 14: India
$ws.Cells.Item(14, 1).Value = "India"
$ws.Cells.Item(14, 2).Value = 118222
$ws.Cells.Item(14, 3).Value = 6194
$ws.Cells.Item(14, 4).Value = 48540
$ws.Cells.Item(14, 5).Value = 66098
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(14, 7).Value = 150
$ws.Cells.Item(14, 8).Value = 3584

# Row 17: Canada
$ws.Cells.Item(17, 1).Value = "Canada"
$ws.Cells.Item(17, 2).Value = 81277
$ws.Cells.Item(17, 3).Value = 1135
$ws.Cells.Item(17, 4).Value = 41603
$ws.Cells.Item(17, 5).Value = 33529
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(17, 7).Value = 114
$ws.Cells.Item(17, 8).Value = 6145

# Row 19: Chile
$ws.Cells.Item(19, 1).Value = "Chile"
$ws.Cells.Item(19, 2).Value = 57581
$ws.Cells.Item(19, 3).Value = 3964
$ws.Cells.Item(19, 4).Value = 23992
$ws.Cells.Item(19, 5).Value = 33000
$ws.Cells.Item(19, 6).Value = 0
$ws.Cells.Item(19, 7).Value = 45
$ws.Cells.Item(19, 8).Value = 589

# Row 20: Mexico
$ws.Cells.Item(20, 1).Value = "Mexico"
$ws.Cells.Item(20, 2).Value = 56594
$ws.Cells.Item(20, 3).Value = 2248
$ws.Cells.Item(20, 4).Value = 38876
$ws.Cells.Item(20, 5).Value = 11628
$ws.Cells.Item(20, 6).Value = 0
$ws.Cells.Item(20, 7).Value = 424
$ws.Cells.Item(20, 8).Value = 6090

# Row 21: Belgica
$ws.Cells.Item(21, 1).Value = "Belgica"
$ws.Cells.Item(21, 2).Value = 56235
$ws.Cells.Item(21, 3).Value = 252
$ws.Cells.Item(21, 4).Value = 14988
$ws.Cells.Item(21, 5).Value = 32061
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(21, 7).Value = 36
$ws.Cells.Item(21, 8).Value = 9186

# Row 28: Suiza
$ws.Cells.Item(28, 1).Value = "Suiza"
$ws.Cells.Item(28, 2).Value = 30694
$ws.Cells.Item(28, 3).Value = 36
$ws.Cells.Item(28, 4).Value = 27800
$ws.Cells.Item(28, 5).Value = 996
$ws.Cells.Item(28, 6).Value = 0
$ws.Cells.Item(28, 7).Value = 6
$ws.Cells.Item(28, 8).Value = 1898

# Row 32: Emiratos Arabes Unidos
$ws.Cells.Item(32, 1).Value = "Emiratos Arabes Unidos"
$ws.Cells.Item(32, 2).Value = 26898
$ws.Cells.Item(32, 3).Value = 894
$ws.Cells.Item(32, 4).Value = 12755
$ws.Cells.Item(32, 5).Value = 13906
$ws.Cells.Item(32, 6).Value = 0
$ws.Cells.Item(32, 7).Value = 4
$ws.Cells.Item(32, 8).Value = 237

# Row 44: Egipto
$ws.Cells.Item(44, 1).Value = "Egipto"
$ws.Cells.Item(44, 2).Value = 15003
$ws.Cells.Item(44, 3).Value = 774
$ws.Cells.Item(44, 4).Value = 4217
$ws.Cells.Item(44, 5).Value = 10090
$ws.Cells.Item(44, 6).Value = 0
$ws.Cells.Item(44, 7).Value = 16
$ws.Cells.Item(44, 8).Value = 696

# Row 135: Republica de Africa Central
$ws.Cells.Item(135, 1).Value = "Republica de Africa Central"
$ws.Cells.Item(135, 2).Value = 436
$ws.Cells.Item(135, 3).Value = 18
$ws.Cells.Item(135, 4).Value = 18
$ws.Cells.Item(135, 5).Value = 418
$ws.Cells.Item(135, 6).Value = 0
$ws.Cells.Item(135, 7).Value = 0
$ws.Cells.Item(135, 8).Value = 0

# Row 136: Congo
$ws.Cells.Item(136, 1).Value = "Congo"
$ws.Cells.Item(136, 2).Value = 420
$ws.Cells.Item(136, 3).Value = 0
$ws.Cells.Item(136, 4).Value = 132
$ws.Cells.Item(136, 5).Value = 273
$ws.Cells.Item(136, 6).Value = 0
$ws.Cells.Item(136, 7).Value = 0
$ws.Cells.Item(136, 8).Value = 15

# Row 146: Ruanda
$ws.Cells.Item(146, 1).Value = "Ruanda"
$ws.Cells.Item(146, 2).Value = 320
$ws.Cells.Item(146, 3).Value = 6
$ws.Cells.Item(146, 4).Value = 217
$ws.Cells.Item(146, 5).Value = 103
$ws.Cells.Item(146, 6).Value = 0
$ws.Cells.Item(146, 7).Value = 0
$ws.Cells.Item(146, 8).Value = 0

# Row 153: Yemen
$ws.Cells.Item(153, 1).Value = "Yemen"
$ws.Cells.Item(153, 2).Value = 197
$ws.Cells.Item(153, 3).Value = 13
$ws.Cells.Item(153, 4).Value = 5
$ws.Cells.Item(153, 5).Value = 159
$ws.Cells.Item(153, 6).Value = 0
$ws.Cells.Item(153, 7).Value = 3
$ws.Cells.Item(153, 8).Value = 33

# Row 154: Martinica
$ws.Cells.Item(154, 1).Value = "Martinica"
$ws.Cells.Item(154, 2).Value = 192
$ws.Cells.Item(154, 3).Value = 0
$ws.Cells.Item(154, 4).Value = 91
$ws.Cells.Item(154, 5).Value = 87
$ws.Cells.Item(154, 6).Value = 0
$ws.Cells.Item(154, 7).Value = 0
$ws.Cells.Item(154, 8).Value = 14

# Row 155: Islas Feroe
$ws.Cells.Item(155, 1).Value = "Islas Feroe"
$ws.Cells.Item(155, 2).Value = 187
$ws.Cells.Item(155, 3).Value = 0
$ws.Cells.Item(155, 4).Value = 187
$ws.Cells.Item(155, 5).Value = 0
$ws.Cells.Item(155, 6).Value = 0
$ws.Cells.Item(155, 7).Value = 0
$ws.Cells.Item(155, 8).Value = 0

# Row 158: Gibraltar
$ws.Cells.Item(158, 1).Value = "Gibraltar"
$ws.Cells.Item(158, 2).Value = 151
$ws.Cells.Item(158, 3).Value = 2
$ws.Cells.Item(158, 4).Value = 146
$ws.Cells.Item(158, 5).Value = 5
$ws.Cells.Item(158, 6).Value = 0
$ws.Cells.Item(158, 7).Value = 0
$ws.Cells.Item(158, 8).Value = 0

# Row 199: Santa Lucia
$ws.Cells.Item(199, 1).Value = "Santa Lucia"
$ws.Cells.Item(199, 2).Value = 18
$ws.Cells.Item(199, 3).Value = 0
$ws.Cells.Item(199, 4).Value = 18
$ws.Cells.Item(199, 5).Value = 0
$ws.Cells.Item(199, 6).Value = 0
$ws.Cells.Item(199, 7).Value = 0
$ws.Cells.Item(199, 8).Value = 0

# Row 200: Belice
$ws.Cells.Item(200, 1).Value = "Belice"
$ws.Cells.Item(200, 2).Value = 18
$ws.Cells.Item(200, 3).Value = 0
$ws.Cells.Item(200, 4).Value = 16
$ws.Cells.Item(200, 5).Value = 0
$ws.Cells.Item(200, 6).Value = 0
$ws.Cells.Item(200, 7).Value = 0
$ws.Cells.Item(200, 8).Value = 2

# Row 209: Seychelles
$ws.Cells.Item(209, 1).Value = "Seychelles"
$ws.Cells.Item(209, 2).Value = 11
$ws.Cells.Item(209, 3).Value = 0
$ws.Cells.Item(209, 4).Value = 11
$ws.Cells.Item(209, 5).Value = 0
$ws.Cells.Item(209, 6).Value = 0
$ws.Cells.Item(209, 7).Value = 0
$ws.Cells.Item(209, 8).Value = 0

# Row 210: Groenlandia
$ws.Cells.Item(210, 1).Value = "Groenlandia"
$ws.Cells.Item(210, 2).Value = 11
$ws.Cells.Item(210, 3).Value = 0
$ws.Cells.Item(210, 4).Value = 11
$ws.Cells.Item(210, 5).Value = 0
$ws.Cells.Item(210, 6).Value = 0
$ws.Cells.Item(210, 7).Value = 0
$ws.Cells.Item(210, 8).Value = 0

# Row 214: Bonaire, San Eustaquio y Saba
$ws.Cells.Item(214, 1).Value = "Bonaire, San Eustaquio y Saba"
$ws.Cells.Item(214, 2).Value = 6
$ws.Cells.Item(214, 3).Value = 0
$ws.Cells.Item(214, 4).Value = 6
$ws.Cells.Item(214, 5).Value = 0
$ws.Cells.Item(214, 6).Value = 0
$ws.Cells.Item(214, 7).Value = 0
$ws.Cells.Item(214, 8).Value = 0

# Row 215: Sahara Occidental
$ws.Cells.Item(215, 1).Value = "Sahara Occidental"
$ws.Cells.Item(215, 2).Value = 6
$ws.Cells.Item(215, 3).Value = 0
$ws.Cells.Item(215, 4).Value = 6
$ws.Cells.Item(215, 5).Value = 0
$ws.Cells.Item(215, 6).Value = 0
$ws.Cells.Item(215, 7).Value = 0
$ws.Cells.Item(215, 8).Value = 0

